$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 183, shifting existing rows 183..211 down to 184..212
$ws.Rows.Item(183).Insert()

# Populate the newly inserted row 183 with the new weekly record
$ws.Range("A183").Value = 10
$ws.Range("B183").Value = "Vega Modelo de Temuco"
$ws.Range("C183").Value = "La Araucanía"
$ws.Range("D183").Value = 44637
$ws.Range("D183").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E183").Value = 9
$ws.Range("F183").Value = "Fruta"
$ws.Range("G183").Value = 100102
$ws.Range("H183").Value = "Cítricos"
$ws.Range("I183").Value = 100102006
$ws.Range("J183").Value = "Pomelo"
$ws.Range("K183").Value = "Start Ruby"
$ws.Range("L183").Value = "Primera"
$ws.Range("M183").Value = 75
$ws.Range("N183").Value = 15000
$ws.Range("O183").Value = 15000
$ws.Range("P183").Value = 15000
$ws.Range("Q183").Value = "$/bandeja 15 kilos granel"
$ws.Range("R183").Value = "Región de O'Higgins"
$ws.Range("S183").Value = 1000
$ws.Range("T183").Value = 15
